$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Tgfb2"
$ws.Range("C2").Value2 = "Acvr1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 1.787349
$ws.Range("H2").Value2 = 5.362047
$ws.Range("I2").Value2 = 0.04925512201701282
$ws.Range("J2").Value2 = 0.04925512201701282
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 5.057757666666666
$ws.Range("N2").Value2 = 15.173273
$ws.Range("O2").Value2 = 0.173378811020062
$ws.Range("P2").Value2 = 0.173378811020062
$ws.Range("Q2").Value2 = 9.039978107759
$ws.Range("R2").Value2 = 81.35980296983101
$ws.Range("S2").Value2 = 0.008539794491957759
$ws.Range("T2").Value2 = 0.008539794491957761
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Tgfb2"
$ws.Range("C3").Value2 = "Acvr1"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 1.787349
$ws.Range("H3").Value2 = 5.362047
$ws.Range("I3").Value2 = 0.04925512201701282
$ws.Range("J3").Value2 = 0.04925512201701282
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 14.247411
$ws.Range("N3").Value2 = 42.742233
$ws.Range("O3").Value2 = 0.4883980890531961
$ws.Range("P3").Value2 = 0.4883980890531961
$ws.Range("Q3").Value2 = 25.465095803439
$ws.Range("R3").Value2 = 229.185862230951
$ws.Range("S3").Value2 = 0.02405610746919107
$ws.Range("T3").Value2 = 0.02405610746919107
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Tgfb2"
$ws.Range("C4").Value2 = "Acvr1"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 1.787349
$ws.Range("H4").Value2 = 5.362047
$ws.Range("I4").Value2 = 0.04925512201701282
$ws.Range("J4").Value2 = 0.04925512201701282
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 9.866548666666667
$ws.Range("N4").Value2 = 29.599646
$ws.Range("O4").Value2 = 0.3382230999267418
$ws.Range("P4").Value2 = 0.3382230999267418
$ws.Range("Q4").Value2 = 17.634965892818
$ws.Range("R4").Value2 = 158.714693035362
$ws.Range("S4").Value2 = 0.01665922005586399
$ws.Range("T4").Value2 = 0.01665922005586399
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Tgfb2"
$ws.Range("C5").Value2 = "Acvr1"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 19.46983
$ws.Range("H5").Value2 = 58.40949000000001
$ws.Range("I5").Value2 = 0.5365425847444997
$ws.Range("J5").Value2 = 0.5365425847444997
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 5.057757666666666
$ws.Range("N5").Value2 = 15.173273
$ws.Range("O5").Value2 = 0.173378811020062
$ws.Range("P5").Value2 = 0.173378811020062
$ws.Range("Q5").Value2 = 98.47368195119667
$ws.Range("R5").Value2 = 886.2631375607701
$ws.Range("S5").Value2 = 0.09302511540463219
$ws.Range("T5").Value2 = 0.09302511540463221
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Tgfb2"
$ws.Range("C6").Value2 = "Acvr1"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 19.46983
$ws.Range("H6").Value2 = 58.40949000000001
$ws.Range("I6").Value2 = 0.5365425847444997
$ws.Range("J6").Value2 = 0.5365425847444997
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 14.247411
$ws.Range("N6").Value2 = 42.742233
$ws.Range("O6").Value2 = 0.4883980890531961
$ws.Range("P6").Value2 = 0.4883980890531961
$ws.Range("Q6").Value2 = 277.39467011013
$ws.Range("R6").Value2 = 2496.55203099117
$ws.Range("S6").Value2 = 0.2620463730848762
$ws.Range("T6").Value2 = 0.2620463730848762
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Tgfb2"
$ws.Range("C7").Value2 = "Acvr1"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 19.46983
$ws.Range("H7").Value2 = 58.40949000000001
$ws.Range("I7").Value2 = 0.5365425847444997
$ws.Range("J7").Value2 = 0.5365425847444997
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 9.866548666666667
$ws.Range("N7").Value2 = 29.599646
$ws.Range("O7").Value2 = 0.3382230999267418
$ws.Range("P7").Value2 = 0.3382230999267418
$ws.Range("Q7").Value2 = 192.1000252267267
$ws.Range("R7").Value2 = 1728.90022704054
$ws.Range("S7").Value2 = 0.1814710962549912
$ws.Range("T7").Value2 = 0.1814710962549912
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Tgfb2"
$ws.Range("C8").Value2 = "Acvr1"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 15.03039733333333
$ws.Range("H8").Value2 = 45.091192
$ws.Range("I8").Value2 = 0.4142022932384875
$ws.Range("J8").Value2 = 0.4142022932384875
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 5.057757666666666
$ws.Range("N8").Value2 = 15.173273
$ws.Range("O8").Value2 = 0.173378811020062
$ws.Range("P8").Value2 = 0.173378811020062
$ws.Range("Q8").Value2 = 76.02010734571289
$ws.Range("R8").Value2 = 684.180966111416
$ws.Range("S8").Value2 = 0.07181390112347202
$ws.Range("T8").Value2 = 0.07181390112347202
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Tgfb2"
$ws.Range("C9").Value2 = "Acvr1"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 15.03039733333333
$ws.Range("H9").Value2 = 45.091192
$ws.Range("I9").Value2 = 0.4142022932384875
$ws.Range("J9").Value2 = 0.4142022932384875
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 14.247411
$ws.Range("N9").Value2 = 42.742233
$ws.Range("O9").Value2 = 0.4883980890531961
$ws.Range("P9").Value2 = 0.4883980890531961
$ws.Range("Q9").Value2 = 214.144248301304
$ws.Range("R9").Value2 = 1927.298234711736
$ws.Range("S9").Value2 = 0.2022956084991289
$ws.Range("T9").Value2 = 0.2022956084991289
$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Tgfb2"
$ws.Range("C10").Value2 = "Acvr1"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 15.03039733333333
$ws.Range("H10").Value2 = 45.091192
$ws.Range("I10").Value2 = 0.4142022932384875
$ws.Range("J10").Value2 = 0.4142022932384875
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 9.866548666666667
$ws.Range("N10").Value2 = 29.599646
$ws.Range("O10").Value2 = 0.3382230999267418
$ws.Range("P10").Value2 = 0.3382230999267418
$ws.Range("Q10").Value2 = 148.2981467686702
$ws.Range("R10").Value2 = 1334.683320918032
$ws.Range("S10").Value2 = 0.1400927836158866
$ws.Range("T10").Value2 = 0.1400927836158866
